# Update Portfolio Overview: refresh the "holding" sheet with the latest
# positions (prices/shares updated, some tickers dropped, HZU.TO added).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("holding")

# The new data set only needs 10 rows (1 header + 9 holdings), down from 17.
# Clear out the now-unused trailing rows first so the sheet's used range
# shrinks back to A1:F10.
$ws.Range("A11:F17").Clear()

# New holdings table (ticker, shares, cost, category, currency, product)
$data = @(
    @("HUZ.TO",        250,    12.37,  "SILVER", "CAD", "ETF"),
    @("0P00016N6T.TO",  109.17, 36.64,  "IT",     "CAD", "MUTUAL"),
    @("HZU.TO",         90,    23.41,  "SILVER", "CAD", "ETF"),
    @("VDE",            57,    51.82,  "OIL",    "USD", "ETF"),
    @("HGU.TO",         75,    32.54,  "GOLD",   "CAD", "ETF"),
    @("ENB.TO",         80,    36.119999999999997, "OIL", "CAD", "STOCK"),
    @("CASH",           0,     26798,  "CASH",   "CAD", "CASH"),
    @("CASH",           0,     17208,  "CASH",   "USD", "CASH"),
    @("SPCE",           70,    16.940000000000001, "SPACE", "USD", "STOCK")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row = $row + 1
}

# Column A now needs to comfortably fit the longest ticker ("0P00016N6T.TO")
$ws.Columns.Item(1).ColumnWidth = 13.75

# Move the active selection to reflect where the editor left off
$ws.Range("C7").Select()
